# Update "想去人数" (want-to-go count) values in column F
# for the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    3  = 3050
    7  = 1668
    11 = 1
    12 = 1370
    13 = 12
    14 = 513
    16 = 31
    17 = 4
    18 = 75
    23 = 3195
    24 = 391
    25 = 132
    26 = 214
    27 = 9
    29 = 94
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
